$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 6: Year 21 -> 22 (numeric), Number 5 -> 1, Issue Date -> 26.FEBRUARY.22
$ws.Range("B6").Value = 22
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = "26.FEBRUARY.22"

# --- Update row 7: Year 21 -> 22 (numeric), Number 6 -> 2, Issue Date -> 26.FEBRUARY.22
$ws.Range("B7").Value = 22
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = "26.FEBRUARY.22"

# --- Remove the old bulk of invoice rows (8-35) that are no longer needed
$ws.Range("A8:D35").ClearContents()

# Rows 8-37 take on the tighter "no-data" row height used elsewhere after the cleanup
$ws.Rows("8:37").RowHeight = 14

# --- Selection left where the user's edit ended up
$ws.Range("A6:I38").Select()
